$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.626.98"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "3.152.92"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.24"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.75"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.538"
$ws.Range("E8").Value = "  +15.17%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.435"
$ws.Range("E10").Value = "  +4.95%  "
$ws.Range("E11").Value = "  +2.34%  "
$ws.Range("E12").Value = "  +3.00%  "
$ws.Range("D13").Value = "3.694.97"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.88"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("E15").Value = "  +3.80%  "
$ws.Range("D16").Value = "58.666.67"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.26"
$ws.Range("E17").Value = "  +3.31%  "
$ws.Range("D18").Value = "3.162.81"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.00"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.15"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.35"
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.79"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.525"
$ws.Range("E24").Value = "  +3.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.62"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.23"
$ws.Range("E28").Value = "  +12.18%  "
$ws.Range("D29").Value = "0.0₃0864"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.10"
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.05"
$ws.Range("E32").Value = "  +2.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.16"
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.31"
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.58"
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  +4.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.10"
$ws.Range("E38").Value = "  -4.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.69"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("E40").Value = "  +2.12%  "
$ws.Range("D41").Value = "2.634.51"
$ws.Range("E41").Value = "  +4.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.25"
$ws.Range("E42").Value = "  +5.91%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0286"
$ws.Range("E44").Value = "  +6.30%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.707"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "3.195.04"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("E48").Value = "  +13.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.20"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.981"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.24"
$ws.Range("E51").Value = "  +1.80%  "
